$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 473), pushing the
# existing rows 473-531 down to 475-533 and growing the used range from
# A1:R531 to A1:R533.
$ws.Rows("473:474").Insert()

# New record 1 (Primera) for the latest reporting date.
$ws.Range("A473").Value = 7
$ws.Range("B473").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C473").Value = "Ñuble"
$ws.Range("D473").Value = 45212
$ws.Range("E473").Value = 16
$ws.Range("F473").Value = 100114013
$ws.Range("G473").Value = "Zanahoria"
$ws.Range("H473").Value = "Sin especificar"
$ws.Range("I473").Value = "Primera"
$ws.Range("J473").Value = 150
$ws.Range("K473").Value = 6500
$ws.Range("L473").Value = 6500
$ws.Range("M473").Value = 6500
$ws.Range("N473").Value = "$/saco 20 kilos"
$ws.Range("O473").Value = "Región de Ñuble"
$ws.Range("P473").Value = 325
$ws.Range("Q473").Value = 20
$ws.Range("R473").Value = "Hortaliza"

# New record 2 (Segunda) for the same reporting date.
$ws.Range("A474").Value = 7
$ws.Range("B474").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C474").Value = "Ñuble"
$ws.Range("D474").Value = 45212
$ws.Range("E474").Value = 16
$ws.Range("F474").Value = 100114013
$ws.Range("G474").Value = "Zanahoria"
$ws.Range("H474").Value = "Sin especificar"
$ws.Range("I474").Value = "Segunda"
$ws.Range("J474").Value = 150
$ws.Range("K474").Value = 5500
$ws.Range("L474").Value = 5500
$ws.Range("M474").Value = 5500
$ws.Range("N474").Value = "$/saco 20 kilos"
$ws.Range("O474").Value = "Región de Ñuble"
$ws.Range("P474").Value = 275
$ws.Range("Q474").Value = 20
$ws.Range("R474").Value = "Hortaliza"
